# Update workbook:
#  - Column C (row 2..260): change the "Förändrad" date serial from 45184 to 45186
#  - Columns S..Y: add a second argument (display text) to every =HYPERLINK("URL") formula,
#    the display text being the URL's file name without its extension.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$used = $ws.UsedRange
$lastRow = $used.Row + $used.Rows.Count - 1

$colC = 3   # column C
$colFirst = 19  # column S
$colLast = 25   # column Y

for ($r = 2; $r -le $lastRow; $r++) {

    # --- Update "Förändrad" date in column C ---
    $cCell = $ws.Cells.Item($r, $colC)
    if ($cCell.Value2 -eq 45184) {
        $cCell.Value = 45186
    }

    # --- Update HYPERLINK formulas in columns S..Y ---
    for ($c = $colFirst; $c -le $colLast; $c++) {
        $cell = $ws.Cells.Item($r, $c)
        if ($cell.HasFormula) {
            $formula = $cell.Formula
            if ($formula.IndexOf("HYPERLINK(") -ge 0 -and $formula.IndexOf(",") -lt 0) {
                $firstQuote = $formula.IndexOf('"')
                $lastQuote = $formula.LastIndexOf('"')
                if ($firstQuote -ge 0 -and $lastQuote -gt $firstQuote) {
                    $url = $formula.Substring($firstQuote + 1, $lastQuote - $firstQuote - 1)

                    $slashIdx = $url.LastIndexOf('/')
                    if ($slashIdx -ge 0) {
                        $fileName = $url.Substring($slashIdx + 1)
                    } else {
                        $fileName = $url
                    }

                    $dotIdx = $fileName.LastIndexOf('.')
                    if ($dotIdx -gt 0) {
                        $baseName = $fileName.Substring(0, $dotIdx)
                    } else {
                        $baseName = $fileName
                    }

                    $newFormula = '=HYPERLINK("' + $url + '", "' + $baseName + '")'
                    $cell.Formula = $newFormula
                }
            }
        }
    }
}
